$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 167, shifting the existing rows (167-251) down to (169-253)
$ws.Rows("167:168").Insert()

# Resize the table (ListObject) so its range/autoFilter grows along with the sheet data
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F253"))

# Row 167: "ASP.net Page" / "Set page title"
$ws.Range("A167").Value = "'@C#"
$ws.Range("B167").Value = "ASP.net Page"
$ws.Range("C167").Value = "Set page title"
$ws.Range("E167").Value = "// set page title`nvar page = HttpContext.Current.Handler as Page;`npage.Title = ""`${1:This page title works}"";"
$ws.Range("E167").WrapText = $true
$ws.Rows("167").RowHeight = 45

# Row 168: "ASP.net Page" / "Set header metadata"
$ws.Range("A168").Value = "'@C#"
$ws.Range("B168").Value = "ASP.net Page"
$ws.Range("C168").Value = "Set header metadata"
$ws.Range("E168").Value = "// set MetaDescription - page variable must be declared before`nvar metaDescription = (HtmlMeta)page.FindControl(""metaDescription"");`nmetaDescription.Content = ""`${1:This MetaDescription works}"";"
$ws.Range("E168").WrapText = $true
$ws.Rows("168").RowHeight = 45

# Update the saved view state (scroll position / active selection)
$ws.Application.ActiveWindow.ScrollRow = 145
$ws.Range("E169").Select()
